{"js": "// This document contains a summary-statistics table with a TTG\n// (tissue transglutaminase) testing breakdown. We need to:\n//   1. Rename \"Days from first to last TTG test\" -> \"Years of testing\"\n//      and update its value.\n//   2. Rename \"Total number of negative TTG\" -> \"Total number of\n//      negative TTG in a row\".\n//   3. Merge the \"0\" counts from the (separate, more granular) existing\n//      \"Total number of negative TTG in a row\" breakdown into the\n//      renamed section above, update rows \"1\"..\"5\" to match that other\n//      breakdown's values, and drop the now-redundant rows/table.\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Load every row + first-cell text for every table so we can locate our\n// target table by content instead of a hard-coded index.\nconst allRows = [];\nfor (let i = 0; i < tables.items.length; i++) {\n  const rows = tables.items[i].rows;\n  rows.load(\"items\");\n  allRows.push(rows);\n}\nawait context.sync();\n\nfor (let i = 0; i < allRows.length; i++) {\n  const rows = allRows[i];\n  for (let r = 0; r < rows.items.length; r++) {\n    rows.items[r].cells.load(\"items\");\n  }\n}\nawait context.sync();\n\nfor (let i = 0; i < allRows.length; i++) {\n  const rows = allRows[i];\n  for (let r = 0; r < rows.items.length; r++) {\n    const cells = rows.items[r].cells.items;\n    if (cells.length > 0) {\n      cells[0].body.load(\"text\");\n    }\n  }\n}\nawait context.sync();\n\nlet rows = null;\nfor (let i = 0; i < allRows.length; i++) {\n  const candidateRows = allRows[i];\n  for (let r = 0; r < candidateRows.items.length; r++) {\n    const cells = candidateRows.items[r].cells.items;\n    if (cells.length > 0 && cells[0].body.text.trim() === \"Total number of negative TTG\") {\n      rows = candidateRows;\n      break;\n    }\n  }\n  if (rows) break;\n}\n\nif (!rows) {\n  throw new Error(\"Could not find table containing 'Total number of negative TTG'\");\n}\n\n// Load full text (both columns) for every row of the target table.\nfor (let r = 0; r < rows.items.length; r++) {\n  const cells = rows.items[r].cells.items;\n  for (let c = 0; c < cells.length; c++) {\n    cells[c].body.load(\"text\");\n  }\n}\nawait context.sync();\n\nfunction cellText(r, c) {\n  return rows.items[r].cells.items[c].body.text.trim();\n}\n\n// Replace a cell's text while preserving its existing paragraph/run\n// formatting (style, bold, alignment, ...) by targeting the cell body's\n// content range rather than re-writing the whole paragraph from scratch.\nfunction setCellText(r, c, newText) {\n  const cell = rows.items[r].cells.items[c];\n  const range = cell.body.getRange(\"Content\");\n  range.insertText(newText, \"Replace\");\n}\n\n// 1) \"Days from first to last TTG test\" -> \"Years of testing\"\n//    \"1066 (643.47, 1669.02)\" -> \"2.92 (1.76, 4.57)\"\nfor (let r = 0; r < rows.items.length; r++) {\n  if (cellText(r, 0) === \"Days from first to last TTG test\") {\n    setCellText(r, 0, \"Years of testing\");\n    setCellText(r, 1, \"2.92 (1.76, 4.57)\");\n    break;\n  }\n}\n\n// 2) Find the header row \"Total number of negative TTG\" (exact match,\n//    not the \"...in a row\" variant that already exists further down) and\n//    rename it.\nlet negHeaderIdx = -1;\nfor (let r = 0; r < rows.items.length; r++) {\n  if (cellText(r, 0) === \"Total number of negative TTG\") {\n    negHeaderIdx = r;\n    break;\n  }\n}\nif (negHeaderIdx === -1) {\n  throw new Error(\"Could not find 'Total number of negative TTG' header row\");\n}\nsetCellText(negHeaderIdx, 0, \"Total number of negative TTG in a row\");\n\n// 3) Update the \"0\"..\"5\" data rows directly under that header with the\n//    merged-in-zero values (sourced from the old second breakdown).\nconst newValues = {\n  \"0\": \"211 (9%)\",\n  \"1\": \"711 (32%)\",\n  \"2\": \"526 (24%)\",\n  \"3\": \"368 (17%)\",\n  \"4\": \"226 (10%)\",\n  \"5\": \"120 (5%)\",\n};\n\nlet cursor = negHeaderIdx + 1;\nwhile (cursor < rows.items.length && Object.prototype.hasOwnProperty.call(newValues, cellText(cursor, 0))) {\n  setCellText(cursor, 1, newValues[cellText(cursor, 0)]);\n  cursor++;\n}\n\n// Everything from `cursor` up through the end of the now-redundant\n// duplicate section gets deleted:\n//   a) any leftover plain numeric rows from the original (longer)\n//      breakdown under the renamed header (old \"6\"/\"7\"), then\n//   b) the whole second \"Total number of negative TTG in a row\" section:\n//      its own bold header row plus its own \"0\"..\"7\" numeric rows.\nconst deleteStart = cursor;\nwhile (cursor < rows.items.length && /^\\d+$/.test(cellText(cursor, 0))) {\n  cursor++;\n}\nif (cursor < rows.items.length && cellText(cursor, 0) === \"Total number of negative TTG in a row\") {\n  cursor++;\n  while (cursor < rows.items.length && /^\\d+$/.test(cellText(cursor, 0))) {\n    cursor++;\n  }\n}\nconst deleteEnd = cursor;\n\n// Delete from the end backwards so earlier row indices stay valid.\nfor (let r = deleteEnd - 1; r >= deleteStart; r--) {\n  rows.items[r].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Get-CellText($table, $row, $col) {\n    return $table.Cell($row, $col).Range.Text.TrimEnd([char]7, [char]13, [char]10)\n}\n\n# Locate the target table: the TTG testing-summary table that has a row\n# whose first cell is exactly \"Total number of negative TTG\" (the one\n# without \"in a row\" yet).\n$t = $null\nfor ($ti = 1; $ti -le $d.Tables.Count; $ti++) {\n    $candidate = $d.Tables.Item($ti)\n    for ($r = 1; $r -le $candidate.Rows.Count; $r++) {\n        if ((Get-CellText $candidate $r 1) -eq \"Total number of negative TTG\") {\n            $t = $candidate\n            break\n        }\n    }\n    if ($t -ne $null) { break }\n}\n\nif ($t -eq $null) {\n    throw \"Could not find target table containing 'Total number of negative TTG'\"\n}\n\n# 1) \"Days from first to last TTG test\" -> \"Years of testing\"\n#    \"1066 (643.47, 1669.02)\" -> \"2.92 (1.76, 4.57)\"\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    if ((Get-CellText $t $r 1) -eq \"Days from first to last TTG test\") {\n        $t.Cell($r, 1).Range.Text = \"Years of testing\"\n        $t.Cell($r, 2).Range.Text = \"2.92 (1.76, 4.57)\"\n        break\n    }\n}\n\n# 2) Find the header row \"Total number of negative TTG\" (exact match) and\n#    rename it to \"Total number of negative TTG in a row\".\n$negHeaderRow = -1\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    if ((Get-CellText $t $r 1) -eq \"Total number of negative TTG\") {\n        $negHeaderRow = $r\n        break\n    }\n}\nif ($negHeaderRow -eq -1) {\n    throw \"Could not find 'Total number of negative TTG' header row\"\n}\n$t.Cell($negHeaderRow, 1).Range.Text = \"Total number of negative TTG in a row\"\n\n# 3) Update the \"0\"..\"5\" data rows immediately below that header with the\n#    merged-in-zero values.\n$newValues = @{\n    \"0\" = \"211 (9%)\";\n    \"1\" = \"711 (32%)\";\n    \"2\" = \"526 (24%)\";\n    \"3\" = \"368 (17%)\";\n    \"4\" = \"226 (10%)\";\n    \"5\" = \"120 (5%)\";\n}\n\n$cursor = $negHeaderRow + 1\nwhile ($cursor -le $t.Rows.Count) {\n    $label = Get-CellText $t $cursor 1\n    if ($newValues.ContainsKey($label)) {\n        $t.Cell($cursor, 2).Range.Text = $newValues[$label]\n        $cursor = $cursor + 1\n    } else {\n        break\n    }\n}\n\n# `$cursor` now addresses the first row to delete: the old \"6\"/\"7\" rows of\n# this table, plus the entire duplicate \"Total number of negative TTG in a\n# row\" table that followed it (its own header row and \"0\"..\"7\" rows).\n# Delete from the last row backwards so earlier indices stay valid.\nfor ($r = $t.Rows.Count; $r -ge $cursor; $r--) {\n    $t.Rows.Item($r).Delete()\n}\n"}
